$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the oldest data point (old row 2) - this shifts every subsequent
# row up by one, matching rows 3..19 -> 2..18 for columns A (date), B (y_0),
# C (y_0_forecast) and D (y_1). This also updates the sheet dimension
# automatically from A1:E19 to A1:E18.
$ws.Rows.Item(2).Delete()

# Column E (y_1_forecast) holds freshly recomputed forecast values (bugfix
# in the naive forecaster), so clear the first four rows (no forecast yet)
# and overwrite the remaining rows with the corrected numbers.
$ws.Range("E2").Value = $null
$ws.Range("E3").Value = $null
$ws.Range("E4").Value = $null
$ws.Range("E5").Value = $null

$ws.Range("E6").Value  = -0.3496173419443749
$ws.Range("E7").Value  = 0.05500386022236903
$ws.Range("E8").Value  = 0.07916875696107883
$ws.Range("E9").Value  = 0.1656566557188155
$ws.Range("E10").Value = 0.2044493994367125
$ws.Range("E11").Value = 0.2685680645708288
$ws.Range("E12").Value = 0.2977174885593792
$ws.Range("E13").Value = 0.0103609600907939
$ws.Range("E14").Value = -1.407243743159736
$ws.Range("E15").Value = 0.2048390592685578
$ws.Range("E16").Value = 0.2152263639657814
$ws.Range("E17").Value = -0.06071040501895997
$ws.Range("E18").Value = 0.160714157635633
